# Updates cryptos list (price/volume columns, and two row swaps) to match
# the latest scrape. Price-like text values that Excel would otherwise
# auto-parse as numbers (e.g. "1.00", "0.978") are forced to remain plain
# text by temporarily applying a Text number format before assigning the
# value, then resetting the cell style back to "Normal" so no residual
# formatting differences are introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.110.42'
$ws.Range('E2').Value = '  -6.02%  '
$ws.Range('D3').Value = '2.452.08'
$ws.Range('E3').Value = '  -8.59%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '538.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.78'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('E8').Value = '  -2.73%  '
$ws.Range('D9').Value = '2.467.76'
$ws.Range('E9').Value = '  -8.16%  '
$ws.Range('E10').Value = '  -5.72%  '
$ws.Range('E11').Value = '  -1.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.42'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.22%  '
$ws.Range('E13').Value = '  -4.27%  '
$ws.Range('D14').Value = '2.890.95'
$ws.Range('E14').Value = '  -8.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '24.08'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -8.17%  '
$ws.Range('D16').Value = '59.044.73'
$ws.Range('E16').Value = '  -5.99%  '
$ws.Range('E17').Value = '  -5.59%  '
$ws.Range('D18').Value = '2.502.74'
$ws.Range('E18').Value = '  -6.74%  '
$ws.Range('E19').Value = '  -5.36%  '
$ws.Range('E20').Value = '  -5.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.17'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.967'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.34%  '
$ws.Range('E23').Value = '  -7.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.75'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.08%  '
$ws.Range('E25').Value = '  -11.71%  '
$ws.Range('E26').Value = '  -5.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.978'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.74'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.98%  '
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.80'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.08%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -8.93%  '
$ws.Range('E31').Value = '  -5.43%  '
$ws.Range('E32').Value = '  -9.21%  '
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '157.69'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.09%  '
$ws.Range('E35').Value = '  -3.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.63'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.21%  '
$ws.Range('E37').Value = '  -8.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.71'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '313.96'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.67%  '
$ws.Range('E40').Value = '  -6.29%  '
$ws.Range('E41').Value = '  -5.28%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.71'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.90%  '
$ws.Range('B43').Value = 'SuiNetwork'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.831'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -10.35%  '
$ws.Range('E44').Value = '  -0.34%  '
$ws.Range('E45').Value = '  -3.64%  '
$ws.Range('E46').Value = '  -2.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0529'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.09%  '
$ws.Range('E48').Value = '  -3.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.13'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.06%  '
$ws.Range('E50').Value = '  -4.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.40'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -8.84%  '
